$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "MetricCollector"

# Populate the header row
$ws.Range("A1").Value = "software_environment"
$ws.Range("B1").Value = "repository"
$ws.Range("C1").Value = "inputs"
$ws.Range("D1").Value = "outputs"
$ws.Range("E1").Value = "id"
$ws.Range("F1").Value = "name"
$ws.Range("G1").Value = "description"
